$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 33: update content in place (Z11_B04_P01 -> Z11_B03_P02)
$ws.Cells.Item(33,1).Value2 = "Z11_B03_P02"
$ws.Cells.Item(33,2).Value2 = "Z11_B03"
$ws.Cells.Item(33,3).Value2 = "Lebensqualität, Gesundheit und Wohlbefinden im urbanen Raum erhalten und steigern"
$ws.Cells.Item(33,4).Value2 = "XXXLebensqualität, Gesundheit und Wohlbefinden im urbanen Raum erhalten und steigern"

# 2. Insert two new rows right after row 36 so the data that used to sit at
#    row 36 (Z12_B03_P01) moves to row 37, and a brand new row 38 is created.
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(38).Insert()

# Copy the formatting from row 36 onto the two freshly inserted rows so the
# new cells keep the same style (border/fill/font/wrap) as the rest of the
# table.
$ws.Range("A36:D36").Copy()
$ws.Range("A37:D38").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Row 36: update content in place (Z12_B03_P01 -> Z12_B02_P02_Ib01)
$ws.Cells.Item(36,1).Value2 = "Z12_B02_P02_Ib01"
$ws.Cells.Item(36,2).Value2 = "Z12_B02"
$ws.Cells.Item(36,3).Value2 = "Werte und Produktionsmuster reflektieren"
$ws.Cells.Item(36,4).Value2 = "XXXWerte und Produktionsmuster reflektieren"

# 4. Row 37: restore the data that previously lived in row 36
$ws.Cells.Item(37,1).Value2 = "Z12_B03_P01"
$ws.Cells.Item(37,2).Value2 = "Z12_B03"
$ws.Cells.Item(37,3).Value2 = "Vorbildwirkung der öffentlichen Hand für nachhaltige öffentliche Beschaffung verwirklichen"
$ws.Cells.Item(37,4).Value2 = "Giving shape to the public sector’s exemplary role in sustainable procurement"

# 5. Row 38: brand new row (Z12_B04_P01)
$ws.Cells.Item(38,1).Value2 = "Z12_B04_P01"
$ws.Cells.Item(38,2).Value2 = "Z12_B04"
$ws.Cells.Item(38,3).Value2 = "Verschuldung von Verbraucherinnen und Verbrauchern – Überlastung vermeiden"
$ws.Cells.Item(38,4).Value2 = "XXXVerschuldung von Verbraucherinnen und Verbrauchern – Überlastung vermeiden"
